# ETRS & Schedule Notify
#
# Refresh the sample/test data baked into the first data row of each of the
# four sheets (phone numbers, dates, message/record ids) with a newer batch
# of values, without disturbing any existing formatting.
#
# NOTE: several of the new values are numeric- or date-looking strings
# (e.g. "0772763163" with a leading zero, "2024-05-28", "2" ...). Assigning
# them straight to Range.Value would let Excel's input-parsing "helpfully"
# reinterpret them as a number/date (dropping the leading zero, turning the
# date text into a serial number, ...), which is not what the source data
# needs - these columns store plain text. Writing the literal through a
# quoted formula and then collapsing it to a value via Copy/PasteSpecial
# keeps the cell's original number format / style untouched while forcing
# the stored type to stay text.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")

# Sheet1: MobileNumber, Date, Date&Time, Enquiry_Date, Enquiry_PhoneNumber,
# User1_MessageId, User1RecId, Lead_PN, Sales_PN
Set-TextValue $ws1.Range("F2")  "0772763163"
Set-TextValue $ws1.Range("N2")  "2024-05-28"
Set-TextValue $ws1.Range("P2")  "2024-05-31 05:00:00 PM"
Set-TextValue $ws1.Range("AC2") "2024-05-28"
Set-TextValue $ws1.Range("AE2") "0838628274"
Set-TextValue $ws1.Range("AK2") "2"
Set-TextValue $ws1.Range("AN2") "126170"
Set-TextValue $ws1.Range("AT2") "3011713182"
Set-TextValue $ws1.Range("AX2") "2139397248"

# Sheet2: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
Set-TextValue $ws2.Range("F2")  "0772763163"
Set-TextValue $ws2.Range("AE2") "0838628274"
Set-TextValue $ws2.Range("AT2") "3011713182"
Set-TextValue $ws2.Range("AX2") "2139397248"

# Sheet3: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
Set-TextValue $ws3.Range("F2")  "0772763163"
Set-TextValue $ws3.Range("AE2") "0838628274"
Set-TextValue $ws3.Range("AT2") "3011713182"
Set-TextValue $ws3.Range("AX2") "2139397248"

# Sheet4: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
Set-TextValue $ws4.Range("F2")  "0772763163"
Set-TextValue $ws4.Range("AE2") "0838628274"
Set-TextValue $ws4.Range("AT2") "3011713182"
Set-TextValue $ws4.Range("AX2") "2139397248"
